# Correct the volume discount amount figures on the "Order" sheet.
# They were previously shown as negative/parenthesised values, e.g.
# "(US$43.02)" and "(US$58.3)"; the correct data is the positive,
# properly-formatted amount, e.g. "US$43.02" and "US$58.30".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order")

$ws.Range("K2").Value = "US$43.02"
$ws.Range("L2").Value = "US$58.30"

# Reflect the cell that was in focus when the data was corrected.
$ws.Range("L6").Select()
